$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 318 (pushes existing rows 318:438 down to 319:439)
$ws.Rows("318:318").Insert()

# Populate the newly inserted row with the new price-report entry
$ws.Range("A318").Value = 4
$ws.Range("B318").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C318").Value = "Los Lagos"
$ws.Range("D318").Value = 44704
$ws.Range("E318").Value = 10
$ws.Range("F318").Value = 100112006
$ws.Range("G318").Value = "Repollo"
$ws.Range("H318").Value = "Crespo record"
$ws.Range("I318").Value = "Primera"
$ws.Range("J318").Value = 500
$ws.Range("K318").Value = 1800
$ws.Range("L318").Value = 1800
$ws.Range("M318").Value = 1800
$ws.Range("N318").Value = "$/unidad"
$ws.Range("O318").Value = "Región del Maule"
$ws.Range("P318").Value = 1800
$ws.Range("Q318").Value = 1
$ws.Range("R318").Value = "Hortaliza"
